$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Folio No" column (column G)
$ws.Columns.Item(7).Delete()
